$wb = $excel.ActiveWorkbook

# --- Sheet: start_price ---
$ws = $wb.Worksheets.Item("start_price")
$ws.Range("A2").Value = 47.75

# --- Sheet: Linear ---
$ws = $wb.Worksheets.Item("Linear")
$ws.Range("B2").Value = 0.004631031518043279   # mu
$ws.Range("B3").Value = -0.07455679099609663   # B
$ws.Range("B4").Value = 1.670791163804068      # sig2
$ws.Range("B5").Value = "[1.0, 0.23888749825856517, 0.03119686013574433, -0.02333558268855407, -0.02357295989966193, 0.0493153886677503, 0.22195673648124387, 0.3970558603749315, 0.20902764682240677, 0.028465332870464943, -0.024261182953712795, -0.035110967120456665, 0.032970400011226204, 0.20802977418976773, 0.3732669625222403, 0.19715888383940092, 0.021533423029683533, -0.019561299833889198, -0.022410797537782978, 0.011975957094097028]"

# --- Sheet: NonLinear ---
$ws = $wb.Worksheets.Item("NonLinear")
$ws.Range("B3").Value = 0.8914814814814814     # p
$ws.Range("B4").Value = 0.02024924804744265    # mu_0
$ws.Range("B5").Value = -0.007991550674743653  # B_0
$ws.Range("B6").Value = 1.801715192005326      # sig2_0
$ws.Range("B7").Value = 0.07961752136450161    # mu_1
$ws.Range("B8").Value = -0.2411708503571491    # B_1
$ws.Range("B9").Value = 1.551190456185142      # sig2_1
$ws.Range("B10").Value = "[1.0, 0.2360975847098175, 0.0310561805928046, -0.022275566167315707, -0.02355235603381758, 0.04935818765286317, 0.2187478876793816, 0.39589960965547155, 0.2056013637770178, 0.027674233653487448, -0.023916194472734756, -0.03460942504467993, 0.03120649198894752, 0.20578932998815852, 0.3709725492709221, 0.1964538262257991, 0.021233777002962546, -0.020012299028381864, -0.020362794074748853, 0.010867226851637655]"
